$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two innings rows (row 2 and row 3) for runs/balls/fours/sixes had been
# swapped; restore them to match the official match activity / scorecard.
# A leading apostrophe keeps each value stored as text (matching the sheet's
# existing numberStoredAsText convention for these numeric-looking strings),
# and ClearFormats() afterwards strips the quote-prefix cell style so the
# cells keep their original (default) formatting.

$ws.Range("C2").Value = "'8"
$ws.Range("D2").Value = "'12"
$ws.Range("E2").Value = "'0"
$ws.Range("F2").Value = "'0"

$ws.Range("C3").Value = "'10"
$ws.Range("D3").Value = "'8"
$ws.Range("E3").Value = "'1"
$ws.Range("F3").Value = "'1"

$ws.Range("C2:F3").ClearFormats()
